# Regenerate column G ("K" / strike count) values for rows 2..74 on Sheet1.
# These are the newly computed "K" values that replace the previous
# "Strike#" derived figures (see commit message: "use K instead of Strike#").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @(2,2,2,1,0,0,3,2,0,1,1,0,2,1,1,2,2,2,1,0,4,1,0,2,0,1,0,3,2,1,0,2,0,1,0,2,1,2,0,1,0,1,0,2,0,1,0,0,0,2,1,1,1,1,1,1,1,0,1,0,0,0,0,1,0,0,0,0,0,2,3,1,3)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
